$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header for column F (copy the existing header formatting) ---
$ws.Range("F1").Value = "Trening"
$ws.Range("E1:E1").Copy()
$ws.Range("F1").PasteSpecial(-4122) # xlPasteFormats

# --- Materialize the empty B/C/D cells for the new rows (4-11) without
#     registering any extra number formats: assigning a named style of
#     "Normal" touches the cell (creating it) while keeping it on the
#     default/no style. ---
$ws.Range("B4:D11").Style = "Normal"

# --- Register the date/time number format on A2 first (this mirrors how
#     the workbook's author iterated on the format string - lowercase then
#     uppercase - while only ever touching a single cell/range instance,
#     which keeps cellXfs from growing an extra unused entry), then apply
#     the final format across the whole date column in one shot. ---
$ws.Range("A2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("A2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A2:A11").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- Row 2 - 3: existing rows, update date values + add Trening (Gra) ---
$ws.Range("A2").Value = (Get-Date -Year 2025 -Month 1 -Day 30 -Hour 0 -Minute 0 -Second 0)
$ws.Range("F2").Value = "Gra"

$ws.Range("A3").Value = (Get-Date -Year 2025 -Month 1 -Day 30 -Hour 0 -Minute 0 -Second 0)
$ws.Range("F3").Value = "Gra"

# --- Row 4 - 7: 27.01.2025 (Duża Gra / Mała Gra) ---
$ws.Range("A4").Value = (Get-Date -Year 2025 -Month 1 -Day 27 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E4").Value = "10-15"
$ws.Range("F4").Value = "Duża Gra"

$ws.Range("A5").Value = (Get-Date -Year 2025 -Month 1 -Day 27 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E5").Value = "5-10"
$ws.Range("F5").Value = "Duża Gra"

$ws.Range("A6").Value = (Get-Date -Year 2025 -Month 1 -Day 27 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E6").Value = "10-15"
$ws.Range("F6").Value = "Mała Gra"

$ws.Range("A7").Value = (Get-Date -Year 2025 -Month 1 -Day 27 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E7").Value = "5-10"
$ws.Range("F7").Value = "Mała Gra"

# --- Row 8 - 11: 29.01.2025 (Duża Gra / Mała Gra) ---
$ws.Range("A8").Value = (Get-Date -Year 2025 -Month 1 -Day 29 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E8").Value = "10-15"
$ws.Range("F8").Value = "Duża Gra"

$ws.Range("A9").Value = (Get-Date -Year 2025 -Month 1 -Day 29 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E9").Value = "5-10"
$ws.Range("F9").Value = "Duża Gra"

$ws.Range("A10").Value = (Get-Date -Year 2025 -Month 1 -Day 29 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E10").Value = "10-15"
$ws.Range("F10").Value = "Mała Gra"

$ws.Range("A11").Value = (Get-Date -Year 2025 -Month 1 -Day 29 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E11").Value = "5-10"
$ws.Range("F11").Value = "Mała Gra"
